$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E (the "DK" column) ---
# E1 header: copy formatting from D1 (bold/filled/centered header style), then overwrite text
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "DK"

# E2:E4 body values = "0.5" (reuse existing text cell so it stays a shared string, not a number)
$ws.Range("C2").Copy($ws.Range("E2"))
$ws.Range("D3").Copy($ws.Range("E3"))
$ws.Range("C4").Copy($ws.Range("E4"))

# --- New row 5 (the "DK" row) ---
# A5 row label: copy formatting from A4 (bold/filled row-label style), then overwrite text
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "DK"

# B5:D5 body values = "0.5"
$ws.Range("C2").Copy($ws.Range("B5"))
$ws.Range("B3").Copy($ws.Range("C5"))
$ws.Range("C2").Copy($ws.Range("D5"))

# E5 diagonal value = "0.0"
$ws.Range("B2").Copy($ws.Range("E5"))

# Match the recorded selection from the edited file
$null = $ws.Range("H7").Select()
